$d = $word.ActiveDocument

$replacements = @(
    @("115÷7=16, 3", "538÷6=89, 4"),
    @("563÷9=62, 5", "989÷6=164, 5"),
    @("269÷5=53, 4", "710÷5=142, 0"),
    @("794÷6=132, 2", "974÷5=194, 4"),
    @("184÷6=30, 4", "757÷5=151, 2"),
    @("113÷3=37, 2", "601÷5=120, 1"),
    @("902÷2=451, 0", "524÷8=65, 4"),
    @("553÷8=69, 1", "730÷5=146, 0"),
    @("690÷3=230, 0", "942÷4=235, 2"),
    @("474÷7=67, 5", "650÷4=162, 2"),
    @("228÷2=114, 0", "393÷2=196, 1"),
    @("189÷4=47, 1", "545÷9=60, 5"),
    @("713÷5=142, 3", "416÷6=69, 2"),
    @("490÷7=70, 0", "221÷6=36, 5"),
    @("526÷3=175, 1", "718÷9=79, 7"),
    @("281÷9=31, 2", "869÷3=289, 2"),
    @("296÷8=37, 0", "825÷5=165, 0"),
    @("505÷4=126, 1", "903÷2=451, 1"),
    @("876÷2=438, 0", "250÷8=31, 2"),
    @("588÷3=196, 0", "851÷5=170, 1"),
    @("765÷8=95, 5", "189÷4=47, 1"),
    @("653÷6=108, 5", "655÷3=218, 1"),
    @("628÷6=104, 4", "420÷3=140, 0"),
    @("276÷9=30, 6", "850÷6=141, 4"),
    @("285÷8=35, 5", "629÷2=314, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
